$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for "2022-Q3" and push the existing
#    quarters down by one row. The A column is a 0-based running index that
#    is independent of the quarter labels, so it needs to be rewritten too.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Insert a blank row above the current row 2 (shifts everything down).
$summary.Rows.Item(2).Insert()

# Copy the formatting (styles/borders/font) from the row that used to be
# row 2 (now row 3, still has the original per-row style) into the new row 2.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new "2022-Q3" row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 26
$summary.Range("D2").Value = 3.95

# Re-sequence column A (0-based index) for every data row now that there are
# 7 quarters (rows 2-8).
$idx = 0
for ($r = 2; $r -le 8; $r++) {
    $summary.Range("A$r").Value = $idx
    $idx = $idx + 1
}

# ---------------------------------------------------------------------------
# 2. Add a new worksheet "2022-Q3" right after "总计" (i.e. before the sheet
#    that is currently "2022-Q2"), holding the quarter's fund holdings.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$q3.Name = "2022-Q3"

# Header row (row 1) - bold/centered/bordered style matches the other sheets.
$headers = @{B="基金代码"; C="基金名称"; D="基金规模"; E="股票总仓位"; F="仓位占比"; G="持有市值(亿元)"; H="仓位排名"}
foreach ($col in @("B","C","D","E","F","G","H")) {
    $q3.Range("$col" + "1").Value = $headers[$col]
}

$rows = @(
    @{A="0"; B="006604"; C="嘉实消费精选股票A"; D="22.30"; E="86.68"; F="4.31"; G="0.9611"; H="9"}
    @{A="1"; B="006605"; C="嘉实消费精选股票C"; D="16.47"; E="86.68"; F="4.31"; G="0.7099"; H="9"}
    @{A="2"; B="004868"; C="交银施罗德股息优化混合"; D="23.52"; E="86.53"; F="2.98"; G="0.7009"; H="10"}
    @{A="3"; B="011069"; C="工银成长精选混合A"; D="12.40"; E="60.10"; F="2.51"; G="0.3112"; H="6"}
    @{A="4"; B="001140"; C="工银总回报灵活配置混合A"; D="6.10"; E="79.57"; F="3.39"; G="0.2068"; H="8"}
    @{A="5"; B="519125"; C="浦银安盛消费升级混合A"; D="2.17"; E="90.82"; F="8.84"; G="0.1918"; H="1"}
    @{A="6"; B="519115"; C="浦银安盛红利精选混合A"; D="2.45"; E="91.44"; F="5.85"; G="0.1433"; H="3"}
    @{A="7"; B="013341"; C="工银核心机遇混合A"; D="6.17"; E="77.79"; F="2.25"; G="0.1388"; H="9"}
    @{A="8"; B="001320"; C="工银丰盈回报灵活配置混合A"; D="2.81"; E="82.02"; F="4.11"; G="0.1155"; H="8"}
    @{A="9"; B="519176"; C="浦银安盛消费升级混合C"; D="1.07"; E="90.82"; F="8.84"; G="0.0946"; H="1"}
    @{A="10"; B="000763"; C="工银新财富灵活配置混合"; D="2.76"; E="92.61"; F="2.96"; G="0.0817"; H="10"}
    @{A="11"; B="013347"; C="工银丰盈回报灵活配置混合C"; D="1.22"; E="82.02"; F="4.11"; G="0.0501"; H="8"}
    @{A="12"; B="013289"; C="工银食品饮料行业混合A"; D="1.00"; E="93.08"; F="4.97"; G="0.0497"; H="8"}
    @{A="13"; B="013342"; C="工银核心机遇混合C"; D="1.88"; E="77.79"; F="2.25"; G="0.0423"; H="9"}
    @{A="14"; B="011070"; C="工银成长精选混合C"; D="1.68"; E="60.10"; F="2.51"; G="0.0422"; H="6"}
    @{A="15"; B="013290"; C="工银食品饮料行业混合C"; D="0.43"; E="93.08"; F="4.97"; G="0.0214"; H="8"}
    @{A="16"; B="010703"; C="财通智选消费股票A"; D="0.48"; E="93.51"; F="3.06"; G="0.0147"; H="7"}
    @{A="17"; B="006644"; C="弘毅远方消费升级混合A"; D="0.41"; E="83.54"; F="3.55"; G="0.0146"; H="5"}
    @{A="18"; B="010704"; C="财通智选消费股票C"; D="0.46"; E="93.51"; F="3.06"; G="0.0141"; H="7"}
    @{A="19"; B="000649"; C="长城久鑫灵活配置混合"; D="0.50"; E="87.25"; F="2.47"; G="0.0124"; H="9"}
    @{A="20"; B="002512"; C="长城久润混合"; D="0.33"; E="85.19"; F="3.33"; G="0.0110"; H="8"}
    @{A="21"; B="011231"; C="光大保德信锦弘混合A"; D="1.95"; E="26.05"; F="0.52"; G="0.0101"; H="5"}
    @{A="22"; B="014422"; C="弘毅远方消费升级混合C"; D="0.13"; E="83.54"; F="3.55"; G="0.0046"; H="5"}
    @{A="23"; B="011232"; C="光大保德信锦弘混合C"; D="0.82"; E="26.05"; F="0.52"; G="0.0043"; H="5"}
    @{A="24"; B="014029"; C="浦银安盛红利精选混合C"; D="0.07"; E="91.44"; F="5.85"; G="0.0041"; H="3"}
    @{A="25"; B="011477"; C="工银总回报灵活配置混合C"; D="0.05"; E="79.57"; F="3.39"; G="0.0017"; H="8"}
)

$r = 2
foreach ($row in $rows) {
    $q3.Range("A$r").Value = [int]$row.A
    $q3.Range("B$r").NumberFormat = "@"
    $q3.Range("B$r").Value = $row.B
    $q3.Range("C$r").Value = $row.C
    $q3.Range("D$r").NumberFormat = "@"
    $q3.Range("D$r").Value = $row.D
    $q3.Range("E$r").NumberFormat = "@"
    $q3.Range("E$r").Value = $row.E
    $q3.Range("F$r").NumberFormat = "@"
    $q3.Range("F$r").Value = $row.F
    $q3.Range("G$r").NumberFormat = "@"
    $q3.Range("G$r").Value = $row.G
    $q3.Range("H$r").Value = [int]$row.H
    $r = $r + 1
}

# Match the header/index-column style used on the other quarter sheets:
# bold, centered, bordered (same style as "总计"'s header/index cells).
# Direct `.Style =` assignment doesn't stick in this COM host, so copy the
# formatting across via PasteSpecial(xlPasteFormats) instead.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("A2").Copy()
$q3.Range("A2:A27").PasteSpecial(-4122)  # xlPasteFormats

Write-Output "ok"
